$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple text replacements (rows unaffected by row count changes) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text  = "0M"       # was 99.96
$t.Rows.Item(2).Cells.Item(1).Range.Text  = "0M"       # was 0.16
$t.Rows.Item(3).Cells.Item(1).Range.Text  = "0M"       # was 396
$t.Rows.Item(4).Cells.Item(1).Range.Text  = "1190"     # was 588
$t.Rows.Item(5).Cells.Item(1).Range.Text  = "0.00001"  # was 0.00003
# row 6 (0.00247) is unchanged

# --- Delete the two rows that held 0.00009 (row 7) and 0.00008 (row 8) ---
$t.Rows.Item(8).Delete()
$t.Rows.Item(7).Delete()

# After deletion, former row 9 (0.00011) is now row 7 and is unchanged.
# Former row 10 (0.00012) is now row 8 -> becomes 0.00005
$t.Rows.Item(8).Cells.Item(1).Range.Text  = "0.00005"
# Former row 11 (0.00012) is now row 9 -> becomes 0.00020
$t.Rows.Item(9).Cells.Item(1).Range.Text  = "0.00020"
# Former row 12 (0.05583) is now row 10 -> becomes 0.00021
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00021"

# --- Insert two brand-new rows after row 10 (before what is now row 11, "100.0") ---
$refRow  = $t.Rows.Item(11)            # currently "100.0"
$newRow2 = $t.Rows.Add($refRow)        # will hold 0.16130, inserted right before refRow
$newRow2.Cells.Item(1).Range.Text = "0.16130"
$newRow1 = $t.Rows.Add($newRow2)       # will hold 0.00023, inserted right before newRow2
$newRow1.Cells.Item(1).Range.Text = "0.00023"

# Row 13 (100.0) and everything through row 43 (10692.1) remain unchanged.

# --- Collapse the three tab-separated summary rows to single values ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.96"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.16"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "396"
